# Update "想去人数" (interest count) figures in the 展览 sheet and the
# aggregated 全部类型 sheet to match the newly generated gh-pages data
# (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# "展览" (Exhibitions) sheet - column F holds the interest-count values.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value  = 12868
$wsExpo.Range("F5").Value  = 85
$wsExpo.Range("F6").Value  = 69
$wsExpo.Range("F7").Value  = 42
$wsExpo.Range("F10").Value = 12791
$wsExpo.Range("F12").Value = 33
$wsExpo.Range("F13").Value = 8660
$wsExpo.Range("F14").Value = 7663
$wsExpo.Range("F16").Value = 93
$wsExpo.Range("F18").Value = 124
$wsExpo.Range("F19").Value = 978
$wsExpo.Range("F21").Value = 14
$wsExpo.Range("F22").Value = 379
$wsExpo.Range("F23").Value = 184

# "全部类型" (All types) sheet - same events, shifted rows because it
# aggregates every category; update the matching rows too.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 12868
$wsAll.Range("F6").Value  = 85
$wsAll.Range("F7").Value  = 69
$wsAll.Range("F8").Value  = 42
$wsAll.Range("F11").Value = 12791
$wsAll.Range("F13").Value = 33
$wsAll.Range("F14").Value = 8660
$wsAll.Range("F15").Value = 7663
$wsAll.Range("F17").Value = 93
$wsAll.Range("F19").Value = 124
$wsAll.Range("F20").Value = 978
$wsAll.Range("F22").Value = 14
$wsAll.Range("F24").Value = 379
$wsAll.Range("F25").Value = 184

Write-Output "updated interest counts"
